$d = $word.ActiveDocument

# --------------------------------------------------------------------
# Change 1: delete the leading "mil" run (the run right before <page>)
# --------------------------------------------------------------------
$milRng = $d.Content
$milRng.Find.Execute("mil", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$milRng.Delete()

# --------------------------------------------------------------------
# Change 2: "...livres <corr>quinze</corr>," -> "...livres quinze [sols],"
# --------------------------------------------------------------------

# delete the "</corr>" closing-tag run
$closeRng = $d.Content
$closeRng.Find.Execute("</corr>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$closeRng.Delete()

# delete the "<corr>" opening-tag run
$openRng = $d.Content
$openRng.Find.Execute("<corr>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$openRng.Delete()

# the two plain-text runs around "quinze" now merge into one run:
# "Et leur reste sur la premiere annee quatre vintg livres quinze"
$mergedRng = $d.Content
$mergedRng.Find.Execute("Et leur reste sur la premiere annee quatre vintg livres quinze", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# insert " [sols]" right after "quinze" (before the comma), as its own run
$insPoint = $d.Range($mergedRng.End, $mergedRng.End)
$insPoint.InsertAfter(" [sols]")
$newRunRng = $d.Range($mergedRng.End, $mergedRng.End + 7)
$newRunRng.Font.Color = -16777216
